$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 8.439499999999995
    "A10" = -21.96789999999999
    "A12" = -21.39249999999999
    "B12" = 5.368300000000001
    "B17" = 5.311499999999996
    "A18" = -22.13900000000001
    "B26" = 4.370600000000003
    "B27" = 6.719000000000005
    "B28" = 6.529400000000002
    "A37" = -22.0414
    "B37" = 6.086300000000001
    "A55" = -21.9917
    "B65" = 6.051300000000001
    "A68" = -21.4633
    "B73" = 8.9346
    "A77" = -20.25709999999999
    "A78" = -19.80689999999998
    "B84" = 5.690800000000003
    "B85" = 5.660499999999997
    "B93" = 5.812199999999997
    "B95" = 5.611000000000001
    "B98" = 5.073900000000009
    "B99" = 5.707899999999996
    "B101" = 5.722899999999998
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
